$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values that changed (semana 11, 17, 19, 30)
$ws.Range("B12").Value = 470
$ws.Range("B18").Value = 457
$ws.Range("B20").Value = 403
$ws.Range("B31").Value = 382

# Add new rows for semana 31, 32, 33 (epi week 32 of 2025 update)
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 312

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 259

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 2
